# Apply updated odds values to Sheet1, matching the supplied diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 2.38
$ws.Range("U2").Value = 1.8
$ws.Range("V2").Value = 1.95
$ws.Range("W2").Value = 8.5
$ws.Range("AC2").Value = 9
$ws.Range("AG2").Value = 351
$ws.Range("AH2").Value = 9.5
$ws.Range("AI2").Value = 17

# Row 4
$ws.Range("Z4").Value = 13
$ws.Range("AH4").Value = 10
$ws.Range("AW4").Value = 6.5

# Row 5
$ws.Range("M5").Value = 1.11
$ws.Range("O5").Value = 1.44
$ws.Range("P5").Value = 2.63

# Row 6
$ws.Range("M6").Value = 1.08
$ws.Range("O6").Value = 1.36

# Row 14
$ws.Range("O14").Value = 1.1

# Row 15
$ws.Range("M15").Value = 1.03
$ws.Range("O15").Value = 1.14

# Row 16
$ws.Range("M16").Value = 1.03
$ws.Range("O16").Value = 1.18

# Row 17
$ws.Range("M17").Value = 1.04
$ws.Range("O17").Value = 1.25

# Row 30
$ws.Range("M30").Value = 1.04
$ws.Range("O30").Value = 1.25

# Row 31
$ws.Range("M31").Value = 1.02
$ws.Range("O31").Value = 1.13

# Row 32
$ws.Range("M32").Value = 1.05
$ws.Range("O32").Value = 1.37
